# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodo Mora (column E) - update to the new sorted list of periods
$ws.Range("E16").Value = "2011"
$ws.Range("E17").Value = "2012"
$ws.Range("E18").Value = "2101"
$ws.Range("E19").Value = "2102"
$ws.Range("E20").Value = "2103"
$ws.Range("E21").Value = "2104"
$ws.Range("E22").Value = "2105"
$ws.Range("E23").Value = "2106"

# Salario Basico (column F) - 2011 and 2106 swap their values
$ws.Range("F16").Value = 120000
$ws.Range("F23").Value = 84000

# Valor Mora (column G) - all updated to the new amount
$ws.Range("G16").Value = 3000000
$ws.Range("G17").Value = 3000000
$ws.Range("G18").Value = 3000000
$ws.Range("G19").Value = 3000000
$ws.Range("G20").Value = 3000000
$ws.Range("G21").Value = 3000000
$ws.Range("G22").Value = 3000000
$ws.Range("G23").Value = 3000000
